# "Url para registrar pedido"
# Add two new POST endpoints to the API inventory sheet:
#   row 7: /productos/crear     -> "Crea un producto"  (kept with the underlined style
#                                   that was already used for the previous last row group)
#   row 8: /pedidos/registrar   -> "Crea un pedido"    (new entry, default style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "POST"
$ws.Range("C7").Value = "/productos/crear"
$ws.Range("D7").Value = "Crea un producto"

$ws.Range("B8").Value = "POST"
$ws.Range("C8").Value = "/pedidos/registrar"
$ws.Range("D8").Value = "Crea un pedido"

# Row 7 carries an underlined font (matches the target formatting)
$ws.Range("B7:D7").Font.Underline = $true

# Leave the selection on the newly added last cell, like the source workbook
$ws.Range("D8").Select()
